# Update specific values in column E (Sheet1) as per the source data correction.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "E3"  = 16.1065
    "E21" = 16.65560000000001
    "E23" = 16.19759999999998
    "E25" = 17.18340000000001
    "E53" = 16.67090000000001
    "E57" = 16.72570000000001
    "E59" = 16.0283
    "E69" = 17.34950000000002
    "E79" = 18.47150000000003
    "E83" = 16.46639999999999
    "E93" = 17.87230000000002
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
